$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$titles = @(
    'Towards a Semantically Useful Definition of Conformance with a Reference Model',
    'Detecting and Resolving Coupling-Related Infrastructure as Code Based Architecture Smells in Microservice Deployments',
    'Correct-By-Construction Microservices with Model-Driven Engineering',
    'A Security Compliance-by-Design Framework Utilizing Reusable Formal Models',
    'Detection Strategies for Microservice Security Tactics',
    'CALint: a Tool for Enforcing the Clean Architecture’s Dependency Rule in Python',
    'Formalization of the AADL Run-Time Services',
    'Assessing Architecture Conformance to Coupling-Related Infrastructure-as-Code Best Practices: Metrics and Case Studies',
    'Assessing Architecture Conformance to Security-Related Practices in Infrastructure as Code Based Deployments',
    'Conformance assessment of Architectural Design Decisions on API endpoint designs derived from domain models',
    'Architecture violations detection and visualization in the continuous integration pipeline',
    'A Two-Level Approach Based on Model Checking to Support Architecture Conformance Checking',
    'Applying a Multi-platform Architectural Conformance Solution in a Real-world Microservice-based System',
    'ArchPython: architecture conformance checking for Python systems',
    'Enforcing Architectural Security Decisions',
    'Monitoring Behavioral Compliance with Architectural Patterns Based on Complex Event Processing',
    'Formal specification and verification of reusable communication models for distributed systems architecture',
    'On the verification of mission-related properties in software-intensive systems-of-systems architectural design',
    'Detecting Architectural Issues During the Continuous Integration Pipeline',
    'Discovering Architectural Rules in Practice',
    'Formalizing Architectural Rules with Ontologies - An Industrial Evaluation',
    'Automated design tests to check Hibernate design recommendations',
    'Architectural and Behavioral Analysis for Cyber Security',
    'Source-Code Divergence Diagnosis Using Constraints and Cryptography',
    'Architecture conformance analysis using model-based testing: A case study approach',
    'Identifying Software Architecture Erosion Through Code Comments',
    'A Secure Testing Based Approach for Mapreduce Frameworks',
    'Formalizing Reusable Communication Models for Distributed Systems Architecture',
    'An Ontology-Based Approach for Documenting and Validating Architecture Rules',
    'On the Understandability of Semantic Constraints for Behavioral Software Architecture Compliance: A Controlled Experiment',
    'Preventing Erosion in Exception Handling Design Using Static-Architecture Conformance Checking',
    'A Holistic Approach for Managed Evolution of Automotive Software Product Line Architectures',
    'On Adequate Behavior-Based Architecture Conformance Checks',
    'Software Architecture Reconstruction and Compliance Checking A Case Study',
    'The accuracy of dependency analysis in static architecture compliance checking',
    'DCL 2.0: modular and reusable specification of architectural constraints',
    'Architecture Conformance Checking with Description Logics',
    'Architectural design of a LMS with LTSA-conformance',
    'Ensuring and Assessing Architecture Conformance to Microservice Decomposition Patterns',
    'Checking sysML models for co-simulation',
    'Modular Specification of Architectural Constraints',
    'Evaluating an Architecture Conformance Monitoring Solution',
    'Behavior-Based Architecture Reconstruction and Conformance Checking',
    'A Genetic Approach to Architectural Pattern Discovery',
    'Architectural Pattern Definition for Semantically Rich Modular Architectures',
    'Dependency Related Parameters in the Reconstruction of a Layered Software Architecture',
    'Rule Type Based Reasoning On Architecture Violations: A Case Study',
    'Preserving architectural decisions through architectural patterns',
    'A Unified Approach to Automatic Testing of Architectural Constraints',
    'A Unified Approach to Architecture Conformance Checking',
    'Preventing Software Architecture Erosion Through Static Architecture Conformance Checking',
    'Automatic translation of architecture constraint specifications into components',
    'Dependency Types and Subtypes in the Context of Architecture Reconstruction and Compliance Checking',
    'OrionPlanning: Improving Modularization and Checking Consistency on Software Architecture',
    'A recommendation system for repairing violations detected by static architecture conformance checking'
)

$startRow = 12
for ($i = 0; $i -lt $titles.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $titles[$i]
}

$endRow = $startRow + $titles.Count - 1
$ws.Range("A11").Copy() | Out-Null
$ws.Range("A" + $startRow + ":A" + $endRow).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A4").Select() | Out-Null
